$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the "Errata- Last modified on ..." date in the first paragraph
#    from 11/8/2018 to 1/23/2019, reproducing the exact run boundaries that
#    Word leaves behind after a sequence of small in-place edits.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("11/8/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1/23/2019", 2) | Out-Null

# Force run splits to land on "Errata- Last modified on 1/" | "23" | "/201" | "9"
# by toggling (and immediately reverting) a character property on each
# sub-range; Word keeps the run boundary even though the formatting ends up
# identical again.
$seg2 = $d.Range(27, 29)
$seg2.Bold = $true
$seg2.Bold = $false

$seg3 = $d.Range(29, 33)
$seg3.Bold = $true
$seg3.Bold = $false

$seg4 = $d.Range(33, 34)
$seg4.Bold = $true
$seg4.Bold = $false

# ---------------------------------------------------------------------------
# 2. Insert the new "Page 128" errata entry right after the "Autobiography"
#    entry (and before the blank line that separates it from "Page 559").
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(4).Range
$anchor.InsertParagraphAfter() | Out-Null

$p = $d.Paragraphs.Item(5)
$p.Range.InsertParagraphAfter() | Out-Null

$p = $d.Paragraphs.Item(6)
$p.Range.Text = "Page 128"
$p.Range.InsertParagraphAfter() | Out-Null

$p = $d.Paragraphs.Item(7)
$p.Range.Text = "In the first full paragraph on the page, in the last clause of the last sentence in that paragraph-"
$p.Range.InsertParagraphAfter() | Out-Null

$p = $d.Paragraphs.Item(8)
$p.Range.Text = "From: but you do not have permission to read it."
$p.Range.InsertParagraphAfter() | Out-Null

$p = $d.Paragraphs.Item(9)
$p.Range.Text = "To: but you do not have permission to write to it."

# ---------------------------------------------------------------------------
# 3. Styles.xml tweaks that came along with this edit in Word:
#    - Normal style gains an explicit (disabled) kerning setting.
#    - Heading style's keepNext paragraph property is reasserted.
# ---------------------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.Font.Kerning = 0

$heading = $d.Styles.Item("Heading")
$heading.ParagraphFormat.KeepWithNext = $true
